$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (row 6 removed entirely)
$data = @(
    @("1005300703", "ПРЕСТИЖ - ВЕ 2015 ООД", 130, "42.82892,23.19246", 42.82892, 23.19246),
    @("1005560104", "ДИМИТРОВ ТРЕЙД СЪРВИС ООД", 145, "42.66627,23.36354", 42.66627, 23.36354),
    @("1005491601", "ДЕСИ СОФИЯ ЕООД", 575, "42.70893,23.38806", 42.70893, 23.38806),
    @("1005355601", "ПЕЦКА ИВАНОВА ЕТ", 60, "42.73325,23.25145", 42.73325, 23.25145)
)

# Column A holds IDs that look purely numeric ("1005300703"); format the
# cells as Text first so the engine keeps them as string cells (matching
# the source's t="inlineStr") instead of coercing them to numbers.
$ws.Range("A2:A5").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Delete the old row 6 entirely (shifts nothing below it since it's the last row)
$ws.Rows.Item(6).Delete()
